$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 11 previously held the "family disease history" / "his_family" mapping.
# For merge-import support the left (label) column now needs to carry the
# raw field name "his_family" as well (matching column B), and the row is
# highlighted in bold to flag it as a special "merge" field.
$ws.Range("A11").Value = "his_family"
$ws.Range("A11:B11").Font.Bold = $true

# Restore the active cell/selection that was recorded when the workbook was
# last saved.
$ws.Range("J15").Select()

# Configure the page for printing (A4 paper, portrait orientation).
$ws.PageSetup.PaperSize = 9
$ws.PageSetup.Orientation = 1
